# Updated test data for DC, TripCurrent, Voltdrop, BatteryStandby
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the TC reference text (was "NGC-1826/TC-63777") ---
$ws.Range("B4").Value = "NGC-1826/T917 OR TC-63777"

# --- New header cells F1 (Loop) and G1 (Column), matching E1's look/fill ---
$null = $ws.Range("E1").Copy()
$null = $ws.Range("F1:G1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Loop"
$ws.Range("G1").Value = "Column"

# --- New data cells F2:F5 (Built-in Loop A-D) and G2 (2), matching B4's plain bordered look ---
$null = $ws.Range("B4").Copy()
$null = $ws.Range("F2:F5").PasteSpecial(-4122)
$null = $ws.Range("G2").PasteSpecial(-4122)

$ws.Range("F2").Value = "Built-in Loop-A"
$ws.Range("G2").Value = 2
$ws.Range("F3").Value = "Built-in Loop-B"
$ws.Range("F4").Value = "Built-in Loop-C"
$ws.Range("F5").Value = "Built-in Loop-D"

# --- Update selection to match new focus area ---
$null = $ws.Range("F1:G5").Select()

$excel.CutCopyMode = 0
